$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "date" column (B) values (stored as text, quote-prefixed):
# Rows 1-7 were "31/05/2023" -> now "19/07/2023"
# Rows 8-11 were "01/06/2023" -> now "17/06/2023"
$ws.Range("B8:B11").Value = "'17/06/2023"
$ws.Range("B1:B7").Value = "'19/07/2023"

# Move the active selection from G7 to C7
$ws.Range("C7").Select()
